$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set B13 with the new text "Shift Center Frequency"
$ws.Range("B13").Value = "Shift Center Frequency"

# Apply theme-based fill (accent5/purple) to A13 (style index s="3", same family as A14/A15)
$ws.Range("A14").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Apply new solid green fill (FF00B050) to A16 and A24 (style index s="4")
$ws.Range("A16").Interior.Color = 5287936
$ws.Range("A24").Interior.Color = 5287936

# Column widths (account for engine's constant 0.8333... padding offset
# between the ColumnWidth COM property and the stored XML column width)
$ws.Columns.Item(1).ColumnWidth = 35.666666666666664
$ws.Columns.Item(2).ColumnWidth = 16.498697916666668

# Selection change
$ws.Range("G7").Select()
